$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 & 19: coin entries swap places (Avalanche <-> ShibaInu), with updated price/volume values
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "13.29"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.000007784"
$ws.Range("E19").Value = "  +0.51%  "

$ws.Range("D2").Value = "30.973.70"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.956.16"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "243.92"
$ws.Range("E5").Value = "  -1.77%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "0.4866"
$ws.Range("E7").Value = "  -0.45%  "
$ws.Range("D8").Value = "0.2939"
$ws.Range("E8").Value = "  -0.71%  "
$ws.Range("D9").Value = "0.07067"
$ws.Range("E9").Value = "  +3.41%  "
$ws.Range("D10").Value = "19.51"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("D11").Value = "107.89"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "1.948.49"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").Value = "0.07752"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "5.368"
$ws.Range("E14").Value = "  -1.54%  "
$ws.Range("D15").Value = "0.7008"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "278.05"
$ws.Range("E16").Value = "  -3.15%  "
$ws.Range("D17").Value = "30.989.31"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D20").Value = "2.212.06"
$ws.Range("E20").Value = "  -0.45%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "5.490"
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "6.513"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").Value = "9.774"
$ws.Range("E25").Value = "  -2.48%  "
$ws.Range("D26").Value = "168.76"
$ws.Range("E26").Value = "  -1.06%  "
$ws.Range("D27").Value = "19.70"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").Value = "0.1048"
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "1.397"
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("D31").Value = "1.564"
$ws.Range("E31").Value = "  -2.39%  "
$ws.Range("D32").Value = "4.589"
$ws.Range("E32").Value = "  -4.98%  "
$ws.Range("D33").Value = "4.413"
$ws.Range("E33").Value = "  -2.35%  "
$ws.Range("D34").Value = "0.04889"
$ws.Range("E34").Value = "  -4.13%  "
$ws.Range("D35").Value = "0.7528"
$ws.Range("E35").Value = "  -2.80%  "
$ws.Range("D36").Value = "1.166"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").Value = "0.02001"
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("D39").Value = "2.680"
$ws.Range("D40").Value = "6.543"
$ws.Range("E40").Value = "  +0.48%  "
$ws.Range("D41").Value = "78.43"
$ws.Range("E41").Value = "  +7.46%  "
$ws.Range("D42").Value = "2.111"
$ws.Range("E42").Value = "  -1.00%  "
$ws.Range("D43").Value = "0.8966"
$ws.Range("E43").Value = "  +0.80%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("D45").Value = "0.4448"
$ws.Range("E45").Value = "  -1.20%  "
$ws.Range("D46").Value = "7.829"
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "988.05"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("E49").Value = "  -1.25%  "
$ws.Range("D50").Value = "9.255"
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").Value = "35.97"
$ws.Range("E51").Value = "  -0.22%  "
